# Auto-generated Excel COM-interop script to apply Alpha_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 514.61536  # H17
$ws.Cells.Item(17, 10).Value = 571.8182  # J17
$ws.Cells.Item(17, 12).Value = 1715.4546  # L17
$ws.Cells.Item(17, 14).Value = -2051.4546  # N17
$ws.Cells.Item(33, 8).Value = 83649  # H33
$ws.Cells.Item(33, 10).Value = 431.66666  # J33
$ws.Cells.Item(33, 12).Value = 431.66666  # L33
$ws.Cells.Item(33, 14).Value = -889.66666  # N33
$ws.Cells.Item(111, 8).Value = 912  # H111
$ws.Cells.Item(111, 9).Value = 912  # I111
$ws.Cells.Item(111, 11).Value = 2736  # K111
$ws.Cells.Item(111, 13).Value = 331  # M111
$ws.Cells.Item(138, 8).Value = 4695  # H138
$ws.Cells.Item(138, 9).Value = 3371.3333  # I138
$ws.Cells.Item(138, 10).Value = 5916.846  # J138
$ws.Cells.Item(138, 11).Value = 10113.9999  # K138
$ws.Cells.Item(138, 12).Value = 17750.538  # L138
$ws.Cells.Item(138, 13).Value = -4973.999899999999  # M138
$ws.Cells.Item(138, 14).Value = -28030.538  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 889.3333  # H4
$ws.Cells.Item(4, 9).Value = 889.3333  # I4
$ws.Cells.Item(4, 11).Value = 889.3333  # K4
$ws.Cells.Item(4, 13).Value = -773.3333  # M4
$ws.Cells.Item(5, 8).Value = 88.166664  # H5
$ws.Cells.Item(5, 9).Value = 89.8  # I5
$ws.Cells.Item(5, 10).Value = 80  # J5
$ws.Cells.Item(5, 11).Value = 89.8  # K5
$ws.Cells.Item(5, 12).Value = 80  # L5
$ws.Cells.Item(5, 13).Value = 22.2  # M5
$ws.Cells.Item(5, 14).Value = -304  # N5
$ws.Cells.Item(32, 8).Value = 3198.9062  # H32
$ws.Cells.Item(32, 9).Value = 3281  # I32
$ws.Cells.Item(32, 11).Value = 3281  # K32
$ws.Cells.Item(32, 13).Value = -2994  # M32
$ws.Cells.Item(45, 8).Value = 1596.8572  # H45
$ws.Cells.Item(45, 9).Value = 1380.5  # I45
$ws.Cells.Item(45, 11).Value = 1380.5  # K45
$ws.Cells.Item(45, 13).Value = -1003.5  # M45
$ws.Cells.Item(80, 8).Value = 26088  # H80
$ws.Cells.Item(83, 8).Value = 26088  # H83
$ws.Cells.Item(101, 8).Value = 48750  # H101
$ws.Cells.Item(101, 10).Value = 48750  # J101
$ws.Cells.Item(101, 12).Value = 48750  # L101
$ws.Cells.Item(101, 14).Value = -55240  # N101
$ws.Cells.Item(132, 8).Value = 100003680  # H132
$ws.Cells.Item(132, 9).Value = 4463.3335  # I132
$ws.Cells.Item(132, 11).Value = 13390.0005  # K132
$ws.Cells.Item(132, 13).Value = -10860.0005  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 88.166664  # H4
$ws.Cells.Item(4, 9).Value = 89.8  # I4
$ws.Cells.Item(4, 10).Value = 80  # J4
$ws.Cells.Item(4, 11).Value = 89.8  # K4
$ws.Cells.Item(4, 12).Value = 80  # L4
$ws.Cells.Item(4, 13).Value = 25.2  # M4
$ws.Cells.Item(4, 14).Value = -310  # N4
$ws.Cells.Item(20, 8).Value = 2381.1304  # H20
$ws.Cells.Item(20, 9).Value = 1772.8  # I20
$ws.Cells.Item(20, 10).Value = 2849.077  # J20
$ws.Cells.Item(20, 11).Value = 1772.8  # K20
$ws.Cells.Item(20, 12).Value = 2849.077  # L20
$ws.Cells.Item(20, 13).Value = -1525.8  # M20
$ws.Cells.Item(20, 14).Value = -3343.077  # N20

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2320.7727  # H31
$ws.Cells.Item(31, 9).Value = 3097  # I31
$ws.Cells.Item(31, 10).Value = 1877.2142  # J31
$ws.Cells.Item(31, 11).Value = 3097  # K31
$ws.Cells.Item(31, 12).Value = 1877.2142  # L31
$ws.Cells.Item(31, 13).Value = -2802  # M31
$ws.Cells.Item(31, 14).Value = -2467.2142  # N31
$ws.Cells.Item(34, 8).Value = 2320.7727  # H34
$ws.Cells.Item(34, 9).Value = 3097  # I34
$ws.Cells.Item(34, 10).Value = 1877.2142  # J34
$ws.Cells.Item(34, 11).Value = 3097  # K34
$ws.Cells.Item(34, 12).Value = 1877.2142  # L34
$ws.Cells.Item(34, 13).Value = -2895  # M34
$ws.Cells.Item(34, 14).Value = -2281.2142  # N34

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2402.889  # H5
$ws.Cells.Item(5, 10).Value = 3499.25  # J5
$ws.Cells.Item(5, 12).Value = 10497.75  # L5
$ws.Cells.Item(5, 14).Value = -10721.75  # N5
$ws.Cells.Item(23, 8).Value = 1373.92  # H23
$ws.Cells.Item(23, 9).Value = 1176.5555  # I23
$ws.Cells.Item(23, 10).Value = 1881.4286  # J23
$ws.Cells.Item(23, 11).Value = 3529.6665  # K23
$ws.Cells.Item(23, 12).Value = 5644.2858  # L23
$ws.Cells.Item(23, 13).Value = -3294.6665  # M23
$ws.Cells.Item(23, 14).Value = -6114.2858  # N23
$ws.Cells.Item(46, 8).Value = 20004978  # H46
$ws.Cells.Item(46, 10).Value = 6331  # J46
$ws.Cells.Item(46, 12).Value = 18993  # L46
$ws.Cells.Item(46, 14).Value = -19175  # N46
$ws.Cells.Item(99, 8).Value = 0  # H99
$ws.Cells.Item(99, 9).Value = 0  # I99
$ws.Cells.Item(99, 10).Value = 0  # J99
$ws.Cells.Item(99, 11).Value = 0  # K99
$ws.Cells.Item(99, 12).Value = 0  # L99
$ws.Cells.Item(99, 13).ClearContents()  # M99
$ws.Cells.Item(99, 14).ClearContents()  # N99
$ws.Cells.Item(118, 8).Value = 1881.1111  # H118
$ws.Cells.Item(118, 9).Value = 1991.25  # I118
$ws.Cells.Item(118, 10).Value = 1000  # J118
$ws.Cells.Item(118, 11).Value = 5973.75  # K118
$ws.Cells.Item(118, 12).Value = 3000  # L118
$ws.Cells.Item(118, 13).Value = -4730.75  # M118
$ws.Cells.Item(118, 14).Value = -5486  # N118
$ws.Cells.Item(121, 8).Value = 1803068.5  # H121
$ws.Cells.Item(121, 10).Value = 2432655.2  # J121
$ws.Cells.Item(121, 12).Value = 7297965.600000001  # L121
$ws.Cells.Item(121, 14).Value = -7300585.600000001  # N121
$ws.Cells.Item(135, 8).Value = 2402.889  # H135
$ws.Cells.Item(135, 10).Value = 3499.25  # J135
$ws.Cells.Item(135, 12).Value = 31493.25  # L135
$ws.Cells.Item(135, 14).Value = -36563.25  # N135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 3354668.8  # H20
$ws.Cells.Item(20, 10).Value = 32003  # J20
$ws.Cells.Item(20, 12).Value = 32003  # L20
$ws.Cells.Item(20, 14).Value = -32493  # N20
$ws.Cells.Item(70, 8).Value = 7581.3477  # H70
$ws.Cells.Item(70, 10).Value = 6698.7856  # J70
$ws.Cells.Item(70, 12).Value = 6698.7856  # L70
$ws.Cells.Item(70, 14).Value = -7238.7856  # N70
$ws.Cells.Item(73, 8).Value = 7581.3477  # H73
$ws.Cells.Item(73, 10).Value = 6698.7856  # J73
$ws.Cells.Item(73, 12).Value = 6698.7856  # L73
$ws.Cells.Item(73, 14).Value = -8570.785599999999  # N73
$ws.Cells.Item(102, 8).Value = 1902  # H102
$ws.Cells.Item(102, 9).Value = 1882.4  # I102
$ws.Cells.Item(102, 11).Value = 1882.4  # K102
$ws.Cells.Item(102, 13).Value = -260.4000000000001  # M102
$ws.Cells.Item(132, 8).Value = 2966.6667  # H132
$ws.Cells.Item(132, 9).Value = 2900  # I132
$ws.Cells.Item(132, 11).Value = 8700  # K132
$ws.Cells.Item(132, 13).Value = -6170  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2478  # H61
$ws.Cells.Item(61, 9).Value = 2478  # I61
$ws.Cells.Item(61, 11).Value = 2478  # K61
$ws.Cells.Item(61, 13).Value = -2276  # M61
$ws.Cells.Item(100, 8).Value = 3000  # H100
$ws.Cells.Item(100, 9).Value = 3000  # I100
$ws.Cells.Item(100, 11).Value = 3000  # K100
$ws.Cells.Item(100, 13).Value = -2459  # M100
$ws.Cells.Item(113, 8).Value = 2478  # H113
$ws.Cells.Item(113, 9).Value = 2478  # I113
$ws.Cells.Item(113, 11).Value = 2478  # K113
$ws.Cells.Item(113, 13).Value = -308  # M113
$ws.Cells.Item(122, 8).Value = 5216.6  # H122
$ws.Cells.Item(122, 10).Value = 4998  # J122
$ws.Cells.Item(122, 12).Value = 14994  # L122
$ws.Cells.Item(122, 14).Value = -19894  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 19636.5  # H45
$ws.Cells.Item(45, 9).Value = 13233  # I45
$ws.Cells.Item(45, 10).Value = 23478.6  # J45
$ws.Cells.Item(45, 11).Value = 13233  # K45
$ws.Cells.Item(45, 12).Value = 23478.6  # L45
$ws.Cells.Item(45, 13).Value = -12742  # M45
$ws.Cells.Item(45, 14).Value = -24460.6  # N45
$ws.Cells.Item(75, 8).Value = 24354.223  # H75
$ws.Cells.Item(75, 9).Value = 19870  # I75
$ws.Cells.Item(75, 10).Value = 24914.75  # J75
$ws.Cells.Item(75, 11).Value = 19870  # K75
$ws.Cells.Item(75, 12).Value = 24914.75  # L75
$ws.Cells.Item(75, 13).Value = -18934  # M75
$ws.Cells.Item(75, 14).Value = -26786.75  # N75
$ws.Cells.Item(78, 8).Value = 24354.223  # H78
$ws.Cells.Item(78, 9).Value = 19870  # I78
$ws.Cells.Item(78, 10).Value = 24914.75  # J78
$ws.Cells.Item(78, 11).Value = 59610  # K78
$ws.Cells.Item(78, 12).Value = 74744.25  # L78
$ws.Cells.Item(78, 13).Value = -54930  # M78
$ws.Cells.Item(78, 14).Value = -84104.25  # N78
$ws.Cells.Item(113, 8).Value = 826.6667  # H113
$ws.Cells.Item(113, 9).Value = 652.5  # I113
$ws.Cells.Item(113, 10).Value = 1697.5  # J113
$ws.Cells.Item(113, 11).Value = 1957.5  # K113
$ws.Cells.Item(113, 12).Value = 5092.5  # L113
$ws.Cells.Item(113, 13).Value = 212.5  # M113
$ws.Cells.Item(113, 14).Value = -9432.5  # N113
$ws.Cells.Item(132, 8).Value = 4610.875  # H132
$ws.Cells.Item(132, 9).Value = 4610.875  # I132
$ws.Cells.Item(132, 11).Value = 13832.625  # K132
$ws.Cells.Item(132, 13).Value = -11302.625  # M132
